$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.305.43"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "2.655.14"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.25"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.79"
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "2.680.89"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.53"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.106"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "3.110.54"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "59.186.54"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.46"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").Value = "2.670.92"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.63"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "347.01"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.56"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.07"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.428"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").Value = "2.767.14"
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("D29").Value = "0.0₃0835"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  +0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.63"
$ws.Range("E31").Value = "  +9.28%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.10"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "149.42"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +16.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.06"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.874"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.66"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "284.35"
$ws.Range("E44").Value = "  -2.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1000"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.991"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.74"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0234"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.993.92"
$ws.Range("E51").Value = "  +1.32%  "
